{"js": "// Office.js (Word JavaScript API) edit script.\n// Replaces the 25 \"three-digit \u00f7 one-digit\" practice answers in the\n// table with their updated values, as described by the commit diff.\n//\n// Body is the async (context) => { ... } function body.\n\nconst replacements = [\n  [\"553\u00f77=79, 0\", \"676\u00f79=75, 1\"],\n  [\"888\u00f79=98, 6\", \"892\u00f76=148, 4\"],\n  [\"893\u00f76=148, 5\", \"167\u00f75=33, 2\"],\n  [\"626\u00f76=104, 2\", \"406\u00f73=135, 1\"],\n  [\"445\u00f73=148, 1\", \"106\u00f79=11, 7\"],\n  [\"226\u00f77=32, 2\", \"711\u00f76=118, 3\"],\n  [\"748\u00f72=374, 0\", \"100\u00f74=25, 0\"],\n  [\"716\u00f79=79, 5\", \"410\u00f73=136, 2\"],\n  [\"160\u00f78=20, 0\", \"586\u00f78=73, 2\"],\n  [\"978\u00f78=122, 2\", \"721\u00f76=120, 1\"],\n  [\"753\u00f78=94, 1\", \"677\u00f73=225, 2\"],\n  [\"354\u00f72=177, 0\", \"509\u00f74=127, 1\"],\n  [\"892\u00f73=297, 1\", \"823\u00f72=411, 1\"],\n  [\"588\u00f78=73, 4\", \"237\u00f73=79, 0\"],\n  [\"791\u00f75=158, 1\", \"157\u00f79=17, 4\"],\n  [\"699\u00f74=174, 3\", \"611\u00f74=152, 3\"],\n  [\"778\u00f79=86, 4\", \"812\u00f78=101, 4\"],\n  [\"804\u00f74=201, 0\", \"650\u00f75=130, 0\"],\n  [\"172\u00f78=21, 4\", \"976\u00f79=108, 4\"],\n  [\"767\u00f76=127, 5\", \"577\u00f79=64, 1\"],\n  [\"453\u00f77=64, 5\", \"141\u00f75=28, 1\"],\n  [\"949\u00f78=118, 5\", \"485\u00f73=161, 2\"],\n  [\"850\u00f76=141, 4\", \"188\u00f79=20, 8\"],\n  [\"643\u00f72=321, 1\", \"429\u00f78=53, 5\"],\n  [\"314\u00f78=39, 2\", \"479\u00f75=95, 4\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText);\n  }\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Replaces the 25 \"three-digit \u00f7 one-digit\" practice answers in the\n# table with their updated values, as described by the commit diff.\n#\n# $word.ActiveDocument is already open as $d below.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"553\u00f77=79, 0\";   New = \"676\u00f79=75, 1\" },\n    @{ Old = \"888\u00f79=98, 6\";   New = \"892\u00f76=148, 4\" },\n    @{ Old = \"893\u00f76=148, 5\";  New = \"167\u00f75=33, 2\" },\n    @{ Old = \"626\u00f76=104, 2\";  New = \"406\u00f73=135, 1\" },\n    @{ Old = \"445\u00f73=148, 1\";  New = \"106\u00f79=11, 7\" },\n    @{ Old = \"226\u00f77=32, 2\";   New = \"711\u00f76=118, 3\" },\n    @{ Old = \"748\u00f72=374, 0\";  New = \"100\u00f74=25, 0\" },\n    @{ Old = \"716\u00f79=79, 5\";   New = \"410\u00f73=136, 2\" },\n    @{ Old = \"160\u00f78=20, 0\";   New = \"586\u00f78=73, 2\" },\n    @{ Old = \"978\u00f78=122, 2\";  New = \"721\u00f76=120, 1\" },\n    @{ Old = \"753\u00f78=94, 1\";   New = \"677\u00f73=225, 2\" },\n    @{ Old = \"354\u00f72=177, 0\";  New = \"509\u00f74=127, 1\" },\n    @{ Old = \"892\u00f73=297, 1\";  New = \"823\u00f72=411, 1\" },\n    @{ Old = \"588\u00f78=73, 4\";   New = \"237\u00f73=79, 0\" },\n    @{ Old = \"791\u00f75=158, 1\";  New = \"157\u00f79=17, 4\" },\n    @{ Old = \"699\u00f74=174, 3\";  New = \"611\u00f74=152, 3\" },\n    @{ Old = \"778\u00f79=86, 4\";   New = \"812\u00f78=101, 4\" },\n    @{ Old = \"804\u00f74=201, 0\";  New = \"650\u00f75=130, 0\" },\n    @{ Old = \"172\u00f78=21, 4\";   New = \"976\u00f79=108, 4\" },\n    @{ Old = \"767\u00f76=127, 5\";  New = \"577\u00f79=64, 1\" },\n    @{ Old = \"453\u00f77=64, 5\";   New = \"141\u00f75=28, 1\" },\n    @{ Old = \"949\u00f78=118, 5\";  New = \"485\u00f73=161, 2\" },\n    @{ Old = \"850\u00f76=141, 4\";  New = \"188\u00f79=20, 8\" },\n    @{ Old = \"643\u00f72=321, 1\";  New = \"429\u00f78=53, 5\" },\n    @{ Old = \"314\u00f78=39, 2\";   New = \"479\u00f75=95, 4\" }\n)\n\nforeach ($rep in $replacements) {\n    $range = $d.Content\n    $found = $range.Find.Execute($rep.Old, $false, $false, $false, $false, $false, $true, 1, $false, $rep.New, 2)\n    if (-not $found) {\n        Write-Output \"WARNING: could not find text to replace: $($rep.Old)\"\n    }\n}\n"}
